$p = $ppt.ActivePresentation
$s = $p.Slides.Item(6)
$sh = $s.Shapes.Item(2)
$tf = $sh.TextFrame
$tr = $tf.TextRange
$tr.Text = "Come premesso nell’introduzione, alcune caratteristiche delle vetture, hanno più impatto rispetto ad altre, sul risultato finale in termini di emissioni.Ad esempio, un motore con un consumo più elevato rispetto ad un altro produrrà più emissioni  di quest’ultimo.Un analisi del dataset, ci mostra quali siano le caratteristiche più impattanti. Un aumento della cilindrata corrisponde ad un aumento delle emissioni. Così come è proporzionale l’aumento delle emissioni in relazione all’aumento dei consumi, etc."
